$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '283.89'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '1.97%'

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '28.55'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '4.41%'

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.102'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '5.35%'

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.06621'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '3.89%'

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '7.293'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '3.81%'

$ws.Range("B7").Value = 'FTXToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.357'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '2.80%'

$ws.Range("B8").Value = 'MXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.9352'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '5.01%'

$ws.Range("B9").Value = 'WazirX'
$ws.Range("C9").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.1564'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '3.09%'

$ws.Range("B10").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C10").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06054'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '10.24%'

$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07604'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '1.55%'

$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.02884'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '-2.59%'

$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.08936'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-0.30%'

$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.001596'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '0.96%'

$ws.Range("B15").Value = 'CoinExToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.04469'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '1.65%'

$ws.Range("B16").Value = 'One'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0006435'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '1.42%'

$ws.Range("B17").Value = 'TigerCash'
$ws.Range("C17").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.006281'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '1.32%'

$ws.Range("B18").Value = 'LEO'
$ws.Range("C18").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.474'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '-0.11%'

$ws.Range("B19").Value = 'GateToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.380'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '2.47%'

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '2.239'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '0.22%'

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.3192'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '0.70%'

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.1301'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '-3.49%'

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.080'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '4.35%'

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.1516'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '0.79%'

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.001177'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '0.19%'

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.004457'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '3.96%'

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0001247'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '5.92%'

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0001608'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '-9.23%'

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.04167'

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.006643'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '-0.60%'

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1246'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '-10.79%'

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002016'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-2.42%'

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.01154'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '3.79%'

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005499'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '-0.84%'

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.01299'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-29.62%'
